$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the first two data rows (2008年 and 2009年), shifting remaining rows up.
$ws.Range("A2:A3").EntireRow.Delete()
